$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the Ost/Nord coordinate values in row 3 to whole numbers
$ws.Range("Q3").Value = 571908
$ws.Range("R3").Value = 6300255

# Clear the Starttid (Z3) and Sluttid (AB3) cells entirely
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
